# mise a jour des packages
#
# Appends a new control-log entry (row 6) to the "Controles" sheet,
# mirroring the structure of the existing rows (2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "19/01/2026"
$ws.Range("B6").Value = "15:31"
$ws.Range("C6").Value = "15:30"
$ws.Range("D6").Value = "poli"
$ws.Range("E6").Value = "Bangoura"
$ws.Range("F6").Value = "Conforme"
$ws.Range("G6").Value = "Conforme"
$ws.Range("H6").Value = "Conforme"
$ws.Range("I6").Value = "Abris bus"
$ws.Range("J6").Value = "Conforme"
$ws.Range("K6").Value = "ras"
$ws.Range("L6").Value = "casc"
$ws.Range("O6").Value = "Lr"

# "125" must be stored as text (matches the other N/A columns' shared
# strings), not auto-converted to a number. A leading apostrophe forces
# text entry; resetting the style back to "Normal" afterwards drops the
# quote-prefix formatting so the cell ends up with the default style.
$ws.Range("P6").Value = "'125"
$ws.Range("P6").Style = "Normal"

$ws.Range("X6").Value = "beau"
$ws.Range("Y6").Value = 102563
$ws.Range("Z6").Value = "Conforme"
$ws.Range("AA6").Value = "Conforme"
$ws.Range("AB6").Value = "Conforme"
$ws.Range("AC6").Value = "Conforme"
$ws.Range("AD6").Value = "Conforme"
$ws.Range("AE6").Value = "Conforme"
$ws.Range("AF6").Value = "Propre"
$ws.Range("AG6").Value = "ras"
$ws.Range("AH6").Value = "Conforme"
$ws.Range("AI6").Value = "Conforme"
$ws.Range("AJ6").Value = "Conforme"
$ws.Range("AK6").Value = "Propre"
$ws.Range("AL6").Value = "Propre"
$ws.Range("AM6").Value = "Propre"
$ws.Range("AN6").Value = "Propre"
$ws.Range("AO6").Value = "ras"
$ws.Range("AP6").Value = 10
$ws.Range("AQ6").Value = 0
$ws.Range("AR6").Value = "BANGOURA"
